# Auto-generated Excel COM-interop edit script
# Applies the DataDictionary.xlsx changes: renames Phase1->Phase2 title,
# relabels duplicate 'Table Name: Library' headers to Member/Rating,
# and appends new Rating/Review/Wishlist table sections (rows 74-92).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells ---
$ws.Cells.Item(2, 1).Value = "DATA DICTIONARY - Iteration Phase 2"
$ws.Cells.Item(55, 1).Value = "Table Name: Member"
$ws.Cells.Item(72, 1).Value = "Table Name: Rating"

# --- Add new table sections (Rating, Review, Wishlist): rows 74-92 ---
# Row 74
$ws.Cells.Item(74, 1).Value = "PK"
$ws.Cells.Item(74, 2).Value = "RatingId"
$ws.Cells.Item(74, 3).Value = "rating key"
$ws.Cells.Item(74, 4).Value = "int"
$ws.Cells.Item(74, 5).Value = "4 bytes"
$ws.Cells.Item(74, 6).Value = "n/a"
# Row 75
$ws.Cells.Item(75, 2).Value = "RatingScore"
$ws.Cells.Item(75, 3).Value = "rating"
$ws.Cells.Item(75, 4).Value = "int"
$ws.Cells.Item(75, 5).Value = "4 bytes"
$ws.Cells.Item(75, 6).Value = "`"Rating`""
# Row 76
$ws.Cells.Item(76, 1).Value = "FK"
$ws.Cells.Item(76, 2).Value = "MemberId"
$ws.Cells.Item(76, 6).Value = "n/a"
# Row 77
$ws.Cells.Item(77, 1).Value = "FK"
$ws.Cells.Item(77, 2).Value = "GameId"
$ws.Cells.Item(77, 6).Value = "n/a"
# Row 79
$ws.Cells.Item(79, 1).Value = "Table Name: Review"
$ws.Cells.Item(79, 1).Font.Bold = $true
$ws.Cells.Item(79, 2).Font.Bold = $true
$ws.Cells.Item(79, 3).Font.Bold = $true
$ws.Cells.Item(79, 4).Font.Bold = $true
$ws.Cells.Item(79, 5).Font.Bold = $true
$ws.Cells.Item(79, 6).Font.Bold = $true
# Row 80
$ws.Cells.Item(80, 1).Value = "PK/FK"
$ws.Cells.Item(80, 1).Font.Bold = $true
$ws.Cells.Item(80, 2).Value = "Field"
$ws.Cells.Item(80, 2).Font.Bold = $true
$ws.Cells.Item(80, 3).Value = "Description"
$ws.Cells.Item(80, 3).Font.Bold = $true
$ws.Cells.Item(80, 4).Value = "Type"
$ws.Cells.Item(80, 4).Font.Bold = $true
$ws.Cells.Item(80, 5).Value = "Size"
$ws.Cells.Item(80, 5).Font.Bold = $true
$ws.Cells.Item(80, 5).HorizontalAlignment = -4108
$ws.Cells.Item(80, 6).Value = "ToolTip"
$ws.Cells.Item(80, 6).Font.Bold = $true
# Row 81
$ws.Cells.Item(81, 1).Value = "PK"
$ws.Cells.Item(81, 2).Value = "ReviewId"
$ws.Cells.Item(81, 3).Value = "review key"
$ws.Cells.Item(81, 4).Value = "int"
$ws.Cells.Item(81, 5).Value = "4 bytes"
$ws.Cells.Item(81, 6).Value = "n/a"
# Row 82
$ws.Cells.Item(82, 2).Value = "Recommended"
$ws.Cells.Item(82, 3).Value = "recommended? (yes/no)"
$ws.Cells.Item(82, 4).Value = "bit"
$ws.Cells.Item(82, 5).Value = "1 byte"
$ws.Cells.Item(82, 6).Value = "`"Recommended`""
# Row 83
$ws.Cells.Item(83, 2).Value = "ReviewText"
$ws.Cells.Item(83, 3).Value = "review text"
$ws.Cells.Item(83, 4).Value = "nvarchar"
$ws.Cells.Item(83, 5).Value = "max"
$ws.Cells.Item(83, 5).HorizontalAlignment = -4108
$ws.Cells.Item(83, 6).Value = "`"Review Text`""
# Row 84
$ws.Cells.Item(84, 2).Value = "Approved"
$ws.Cells.Item(84, 3).Value = "approved flag (hidden)"
$ws.Cells.Item(84, 4).Value = "bit"
$ws.Cells.Item(84, 5).Value = "1 byte"
$ws.Cells.Item(84, 6).Value = "n/a"
# Row 85
$ws.Cells.Item(85, 1).Value = "FK"
$ws.Cells.Item(85, 2).Value = "MemberId"
$ws.Cells.Item(85, 6).Value = "n/a"
# Row 86
$ws.Cells.Item(86, 1).Value = "FK"
$ws.Cells.Item(86, 2).Value = "GameId"
$ws.Cells.Item(86, 6).Value = "n/a"
# Row 88
$ws.Cells.Item(88, 1).Value = "Table Name: Wishlist"
$ws.Cells.Item(88, 1).Font.Bold = $true
$ws.Cells.Item(88, 2).Font.Bold = $true
$ws.Cells.Item(88, 3).Font.Bold = $true
$ws.Cells.Item(88, 4).Font.Bold = $true
$ws.Cells.Item(88, 5).Font.Bold = $true
$ws.Cells.Item(88, 6).Font.Bold = $true
# Row 89
$ws.Cells.Item(89, 1).Value = "PK/FK"
$ws.Cells.Item(89, 1).Font.Bold = $true
$ws.Cells.Item(89, 2).Value = "Field"
$ws.Cells.Item(89, 2).Font.Bold = $true
$ws.Cells.Item(89, 3).Value = "Description"
$ws.Cells.Item(89, 3).Font.Bold = $true
$ws.Cells.Item(89, 4).Value = "Type"
$ws.Cells.Item(89, 4).Font.Bold = $true
$ws.Cells.Item(89, 5).Value = "Size"
$ws.Cells.Item(89, 5).Font.Bold = $true
$ws.Cells.Item(89, 5).HorizontalAlignment = -4108
$ws.Cells.Item(89, 6).Value = "ToolTip"
$ws.Cells.Item(89, 6).Font.Bold = $true
# Row 90
$ws.Cells.Item(90, 1).Value = "PK"
$ws.Cells.Item(90, 2).Value = "WishlistId"
$ws.Cells.Item(90, 3).Value = "wishlist key"
$ws.Cells.Item(90, 4).Value = "int"
$ws.Cells.Item(90, 5).Value = "4 bytes"
$ws.Cells.Item(90, 6).Value = "n/a"
# Row 91
$ws.Cells.Item(91, 1).Value = "FK"
$ws.Cells.Item(91, 2).Value = "MemberId"
$ws.Cells.Item(91, 6).Value = "n/a"
# Row 92
$ws.Cells.Item(92, 1).Value = "FK"
$ws.Cells.Item(92, 2).Value = "GameId"
$ws.Cells.Item(92, 6).Value = "n/a"

# --- Update sheet view: scroll to top and select A2 ---
$ws.Range("A2").Select()
